$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 111611146
$ws.Range("B9").Value = 88630
$ws.Range("E9").Value = 4823
$ws.Range("F9").Value = 'Hasselsopp'
$ws.Range("G9").Value = 'Leccinellum pseudoscabrum'
$ws.Range("H9").Value = '(Kallenb.) Mikšík'
$ws.Range("Q9").Value = 663088.0668624006
$ws.Range("R9").Value = 6634684.960451891
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = '2023-08-11'
$ws.Range("Y9").Style = "Normal"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '2023-08-11'
$ws.Range("AA9").Style = "Normal"
$ws.Range("AC9").Value = '1 ex. under ek och hassel.'
$ws.Range("AX9").Value = 'Gillis Aronsson'
# Row 10
$ws.Range("A10").Value = 111611138
$ws.Range("B10").Value = 81796
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 5406
$ws.Range("F10").Value = 'Gulmjölkig storskål'
$ws.Range("G10").Value = 'Peziza succosa'
$ws.Range("H10").Value = 'Berk.'
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = '3'
$ws.Range("I10").Style = "Normal"
$ws.Range("Q10").Value = 663213.3366271106
$ws.Range("R10").Value = 6634830.464506784
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = '2023-08-12'
$ws.Range("Y10").Style = "Normal"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = '2023-08-12'
$ws.Range("AA10").Style = "Normal"
$ws.Range("AC10").Value = '3 ex. på bar jord och i lövförna.'
$ws.Range("AX10").Value = 'Gillis Aronsson, Cajsa Björkén'
# Row 11
$ws.Range("A11").Value = 111611145
$ws.Range("B11").Value = 88630
$ws.Range("E11").Value = 4823
$ws.Range("F11").Value = 'Hasselsopp'
$ws.Range("G11").Value = 'Leccinellum pseudoscabrum'
$ws.Range("H11").Value = '(Kallenb.) Mikšík'
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = '1'
$ws.Range("I11").Style = "Normal"
$ws.Range("Q11").Value = 663143.8264147732
$ws.Range("R11").Value = 6634793.669287071
$ws.Range("AC11").Value = '1 ex. i lövförna under hassel.'
# Row 12
$ws.Range("A12").Value = 111611165
$ws.Range("B12").Value = 84741
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 37
$ws.Range("F12").Value = 'Jättekamskivling'
$ws.Range("G12").Value = 'Amanita ceciliae'
$ws.Range("H12").Value = '(Berk. & Broome) Bas'
# Row 13
$ws.Range("A13").Value = 111611158
$ws.Range("B13").Value = 86021
$ws.Range("E13").Value = 4037
$ws.Range("F13").Value = 'Bolmörtsskivling'
$ws.Range("G13").Value = 'Entoloma sinuatum'
$ws.Range("H13").Value = '(Bull.) P.Kumm.'
$ws.Range("Q13").Value = 663128.0992466732
$ws.Range("R13").Value = 6634761.25188593
$ws.Range("AC13").Value = '1 ex. i lövförna under ek och hassel.'
